$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 16 & 17: ShibaInu / WrappedEther swapped positions, with updated price & volume ---
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$c = $ws.Range("D16"); $c.NumberFormat = "@"; $c.Value = '1.414.21'; $c.ClearFormats()
$ws.Range("E16").Value = "  -8.07%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$c = $ws.Range("D17"); $c.NumberFormat = "@"; $c.Value = '0.00001020'; $c.ClearFormats()
$ws.Range("E17").Value = "  -6.11%  "

# --- All other Price (D) / Volume(1h) (E) updates ---
$c = $ws.Range("D2"); $c.NumberFormat = "@"; $c.Value = '19.913.97'; $c.ClearFormats()
$ws.Range("E2").Value = '  -7.44%  '
$c = $ws.Range("D3"); $c.NumberFormat = "@"; $c.Value = '1.409.14'; $c.ClearFormats()
$ws.Range("E3").Value = '  -7.97%  '
$c = $ws.Range("D4"); $c.NumberFormat = "@"; $c.Value = '1.005'; $c.ClearFormats()
$ws.Range("E4").Value = '  +0.32%  '
$c = $ws.Range("D5"); $c.NumberFormat = "@"; $c.Value = '1.000'; $c.ClearFormats()
$ws.Range("E5").Value = '  -0.18%  '
$c = $ws.Range("D6"); $c.NumberFormat = "@"; $c.Value = '274.79'; $c.ClearFormats()
$ws.Range("E6").Value = '  -4.66%  '
$c = $ws.Range("D7"); $c.NumberFormat = "@"; $c.Value = '0.3659'; $c.ClearFormats()
$ws.Range("E7").Value = '  -5.78%  '
$c = $ws.Range("D8"); $c.NumberFormat = "@"; $c.Value = '0.3109'; $c.ClearFormats()
$ws.Range("E8").Value = '  -2.06%  '
$c = $ws.Range("D9"); $c.NumberFormat = "@"; $c.Value = '39.72'; $c.ClearFormats()
$ws.Range("E9").Value = '  -6.75%  '
$c = $ws.Range("D10"); $c.NumberFormat = "@"; $c.Value = '1.034'; $c.ClearFormats()
$ws.Range("E10").Value = '  -3.28%  '
$c = $ws.Range("D11"); $c.NumberFormat = "@"; $c.Value = '0.06495'; $c.ClearFormats()
$ws.Range("E11").Value = '  -9.18%  '
$c = $ws.Range("D12"); $c.NumberFormat = "@"; $c.Value = '0.9980'; $c.ClearFormats()
$ws.Range("E12").Value = '  -0.37%  '
$c = $ws.Range("D13"); $c.NumberFormat = "@"; $c.Value = '5.494'; $c.ClearFormats()
$ws.Range("E13").Value = '  -4.07%  '
$c = $ws.Range("D14"); $c.NumberFormat = "@"; $c.Value = '17.69'; $c.ClearFormats()
$ws.Range("E14").Value = '  -2.62%  '
$c = $ws.Range("D15"); $c.NumberFormat = "@"; $c.Value = '6.190'; $c.ClearFormats()
$ws.Range("E15").Value = '  -5.27%  '
$c = $ws.Range("D18"); $c.NumberFormat = "@"; $c.Value = '0.05688'; $c.ClearFormats()
$ws.Range("E18").Value = '  -13.89%  '
$c = $ws.Range("D19"); $c.NumberFormat = "@"; $c.Value = '1.000'; $c.ClearFormats()
$ws.Range("E19").Value = '  -0.19%  '
$c = $ws.Range("D20"); $c.NumberFormat = "@"; $c.Value = '70.92'; $c.ClearFormats()
$ws.Range("E20").Value = '  -14.84%  '
$c = $ws.Range("D21"); $c.NumberFormat = "@"; $c.Value = '5.633'; $c.ClearFormats()
$ws.Range("E21").Value = '  -7.70%  '
$c = $ws.Range("D22"); $c.NumberFormat = "@"; $c.Value = '14.70'; $c.ClearFormats()
$ws.Range("E22").Value = '  -4.51%  '
$c = $ws.Range("D23"); $c.NumberFormat = "@"; $c.Value = '11.01'; $c.ClearFormats()
$ws.Range("E23").Value = '  +1.59%  '
$c = $ws.Range("D24"); $c.NumberFormat = "@"; $c.Value = '2.261'; $c.ClearFormats()
$ws.Range("E24").Value = '  -4.53%  '
$c = $ws.Range("D25"); $c.NumberFormat = "@"; $c.Value = '19.989.06'; $c.ClearFormats()
$ws.Range("E25").Value = '  -7.11%  '
$c = $ws.Range("D26"); $c.NumberFormat = "@"; $c.Value = '2.264'; $c.ClearFormats()
$ws.Range("E26").Value = '  -4.58%  '
$c = $ws.Range("D27"); $c.NumberFormat = "@"; $c.Value = '133.49'; $c.ClearFormats()
$ws.Range("E27").Value = '  -10.13%  '
$c = $ws.Range("D28"); $c.NumberFormat = "@"; $c.Value = '17.19'; $c.ClearFormats()
$ws.Range("E28").Value = '  -6.25%  '
$c = $ws.Range("D29"); $c.NumberFormat = "@"; $c.Value = '1.571.22'; $c.ClearFormats()
$ws.Range("E29").Value = '  -8.01%  '
$c = $ws.Range("D30"); $c.NumberFormat = "@"; $c.Value = '109.88'; $c.ClearFormats()
$ws.Range("E30").Value = '  -5.76%  '
$c = $ws.Range("D31"); $c.NumberFormat = "@"; $c.Value = '3.951'; $c.ClearFormats()
$ws.Range("E31").Value = '  -18.24%  '
$c = $ws.Range("D32"); $c.NumberFormat = "@"; $c.Value = '5.305'; $c.ClearFormats()
$ws.Range("E32").Value = '  -12.42%  '
$c = $ws.Range("D33"); $c.NumberFormat = "@"; $c.Value = '0.8215'; $c.ClearFormats()
$ws.Range("E33").Value = '  -14.04%  '
$c = $ws.Range("D34"); $c.NumberFormat = "@"; $c.Value = '0.07690'; $c.ClearFormats()
$ws.Range("E34").Value = '  -3.95%  '
$c = $ws.Range("D35"); $c.NumberFormat = "@"; $c.Value = '8.397'; $c.ClearFormats()
$ws.Range("E35").Value = '  -1.37%  '
$c = $ws.Range("D36"); $c.NumberFormat = "@"; $c.Value = '1.484'; $c.ClearFormats()
$ws.Range("E36").Value = '  -0.31%  '
$c = $ws.Range("D37"); $c.NumberFormat = "@"; $c.Value = '0.05898'; $c.ClearFormats()
$ws.Range("E37").Value = '  +0.20%  '
$c = $ws.Range("D38"); $c.NumberFormat = "@"; $c.Value = '4.902'; $c.ClearFormats()
$ws.Range("E38").Value = '  -5.10%  '
$c = $ws.Range("D39"); $c.NumberFormat = "@"; $c.Value = '0.9998'; $c.ClearFormats()
$ws.Range("E39").Value = '  -0.20%  '
$c = $ws.Range("D40"); $c.NumberFormat = "@"; $c.Value = '0.02071'; $c.ClearFormats()
$ws.Range("E40").Value = '  -5.72%  '
$c = $ws.Range("D41"); $c.NumberFormat = "@"; $c.Value = '10.48'; $c.ClearFormats()
$ws.Range("E41").Value = '  -7.14%  '
$c = $ws.Range("D42"); $c.NumberFormat = "@"; $c.Value = '0.1900'; $c.ClearFormats()
$ws.Range("E42").Value = '  -5.99%  '
$c = $ws.Range("D43"); $c.NumberFormat = "@"; $c.Value = '1.095'; $c.ClearFormats()
$ws.Range("E43").Value = '  -7.21%  '
$c = $ws.Range("D44"); $c.NumberFormat = "@"; $c.Value = '12.41'; $c.ClearFormats()
$ws.Range("E44").Value = '  -5.13%  '
$c = $ws.Range("D45"); $c.NumberFormat = "@"; $c.Value = '0.5304'; $c.ClearFormats()
$ws.Range("E45").Value = '  -7.75%  '
$c = $ws.Range("D46"); $c.NumberFormat = "@"; $c.Value = '3.529'; $c.ClearFormats()
$ws.Range("E46").Value = '  -5.04%  '
$c = $ws.Range("D47"); $c.NumberFormat = "@"; $c.Value = '0.5200'; $c.ClearFormats()
$ws.Range("E47").Value = '  -6.44%  '
$c = $ws.Range("D48"); $c.NumberFormat = "@"; $c.Value = '115.79'; $c.ClearFormats()
$ws.Range("E48").Value = '  +0.25%  '
$c = $ws.Range("D49"); $c.NumberFormat = "@"; $c.Value = '1.768'; $c.ClearFormats()
$ws.Range("E49").Value = '  -6.67%  '
$c = $ws.Range("D50"); $c.NumberFormat = "@"; $c.Value = '1.038'; $c.ClearFormats()
$ws.Range("E50").Value = '  -10.45%  '
$c = $ws.Range("D51"); $c.NumberFormat = "@"; $c.Value = '0.9983'; $c.ClearFormats()
$ws.Range("E51").Value = '  -0.44%  '
